$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("toBean")

$ws.Range("C1").Value = "list#key?toMap=key"
$ws.Range("D1").Value = "list#value?toMap=value&type=string"
